# Update the cached "datetimeFigureOut" date placeholder text that appears
# on the slide master and on every slide layout (the footer/date placeholder,
# PlaceholderFormat.Type = 16 -> ppPlaceholderDate) from
# "2018/11/19 Monday" to "2018/12/12 Wednesday".

$p = $ppt.ActivePresentation
$m = $p.SlideMaster

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq 16) {
            if ($sh.HasTextFrame -eq -1) {
                $sh.TextFrame.TextRange.Text = "2018/12/12 Wednesday"
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $m.Shapes

# Every slide layout belonging to the master
for ($j = 1; $j -le $m.CustomLayouts.Count; $j++) {
    $layout = $m.CustomLayouts.Item($j)
    Update-DatePlaceholder $layout.Shapes
}
